# lezione4 edits: inserts several new paragraphs of lecture notes, adds a
# lastRenderedPageBreak marker, and restructures the EX-facebook / OSINT
# paragraphs (moving the stray "_GoBack" bookmark).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParagraphsAfter($anchorText, [string[]]$xmlFragments) {
    # Finds the paragraph whose Range text starts with $anchorText and
    # inserts one new (empty) paragraph per fragment right after it, then
    # fills each with the corresponding literal OOXML.
    $found = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($anchorText)) {
            $found = $p
            break
        }
    }
    if ($found -eq $null) {
        throw "anchor paragraph not found: $anchorText"
    }
    $r = $found.Range
    $r.Collapse(0) | Out-Null
    foreach ($frag in $xmlFragments) {
        $r.InsertParagraphAfter() | Out-Null
        $next = $found.Next()
        $next.Range.InsertXML($frag) | Out-Null
        $found = $next
        $r = $found.Range
        $r.Collapse(0) | Out-Null
    }
}

# --- 1. After "Mi sono svegliato tardi..." insert four new paragraphs ---
$p_a = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Durante questi venti minuti parla principalmente dell’introduzione con definizione di </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>alibi(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>non-presenza, presenza fisica o scusanti che attestino di non avere fatto un reato, reato in un luogo diverso), l’alibi informatico è teoricamente facilmente contraffabile, è molto semplice fare un alibi informatico.</w:t></w:r></w:p>'
$p_b = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Tecniche di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>antiforensics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$p_c = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Riscontri oggettivi per cui l’indagato si trovava in un posto diverso rispetto al delitto, è una </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>giustificazione ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> se è totalmente attestato l’indagato viene totalmente scagionato dall’investigazione.</w:t></w:r></w:p>'
$p_d = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Trovarsi in un altro luogo, sia virtuale in modo da non aver compiuto il diritto, che fisico per non essere nel luogo del </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>delitto(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>sempre se ha un luogo ben definito)</w:t></w:r></w:p>'

Insert-ParagraphsAfter "Mi sono svegliato tardi" @($p_a, $p_b, $p_c, $p_d)

# --- 2. After "Ci sono molti modi di simulare un alibi..." insert one paragraph ---
$p_e = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Cosa ci da la sicurezza che colui che ha fatto le azioni che costruiscono l’alibi sia effettivamente il proprietario delle </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>credenziali(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>complici, bot, script)</w:t></w:r></w:p>'

Insert-ParagraphsAfter "Ci sono molti modi di simulare un alibi" @($p_e)

# --- 3. After "Per smontare questi alibi..." insert two paragraphs ---
$p_f = '<w:p ' + $wNs + '><w:r><w:t>Un perito o un consulente tecnico deve fare sempre tutto in buona fede, se sbaglia ci può stare, ma se mente di propria volontà o perché costretto è penalmente conseguibile</w:t></w:r></w:p>'
$p_g = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">È possibile risalire agli spostamenti del soggetto anche da celle telefoniche, non solo dal </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, anche se è abbastanza più difficile</w:t></w:r></w:p>'

Insert-ParagraphsAfter "Per smontare questi alibi" @($p_f, $p_g)

# --- 4. After "sindrome dell'abitacolo..." insert one paragraph ---
$p_h = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Altri reati su internet sono </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>infiniti(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>revenge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>porn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, rilascio di informazioni confidenziali, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…)</w:t></w:r></w:p>'

Insert-ParagraphsAfter "sindrome dell" @($p_h)

Write-Host "done"
